# Add "Gameweeks import" columns to the Challenges sheet:
#   R: "Show Statistics Continuously" (header) / "true" (data row, as text)
#   S: "Gameweek" (header) / 1 (data row, numeric)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

$ws.Range("R1").Value = "Show Statistics Continuously"
$ws.Range("S1").Value = "Gameweek"

# Leading apostrophe forces this to be stored as the literal text "true"
# instead of being auto-coerced to the Boolean TRUE; resetting the style
# back to Normal drops the quote-prefix formatting afterwards.
$ws.Range("R2").Value = "'true"
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").Value = 1
